$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a "clean" numeric-looking string (e.g. "480.20")
# must be forced to Text format first, otherwise Excel auto-converts the
# assigned value to a number and silently drops the significant trailing
# zeros / formatting (e.g. "480.20" -> 480.2, "1.00" -> 1).
$textCells = @("D4", "D5", "D6", "D7", "D9", "D10", "D11", "D12", "D13", "D15", "D18", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D32", "D33", "D35", "D36", "D39", "D41", "D42", "D43", "D44", "D45", "D47", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values (prices, volumes, names and links).
$ws.Range("D2").Value = "68.311.58"
$ws.Range("E2").Value = "  +1.97%  "
$ws.Range("D3").Value = "3.901.55"
$ws.Range("E3").Value = "  +1.52%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("D5").Value = "480.20"
$ws.Range("E5").Value = "  +3.07%  "
$ws.Range("D6").Value = "144.76"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").Value = "0.620"
$ws.Range("E7").Value = "  -1.71%  "
$ws.Range("D9").Value = "0.723"
$ws.Range("E9").Value = "  -2.69%  "
$ws.Range("D10").Value = "0.166"
$ws.Range("E10").Value = "  +7.93%  "
$ws.Range("D11").Value = "0.0000351"
$ws.Range("E11").Value = "  +14.10%  "
$ws.Range("D12").Value = "42.64"
$ws.Range("E12").Value = "  -1.46%  "
$ws.Range("D13").Value = "10.65"
$ws.Range("E13").Value = "  +2.26%  "
$ws.Range("D14").Value = "4.535.58"
$ws.Range("E14").Value = "  +1.04%  "
$ws.Range("D15").Value = "14.59"
$ws.Range("E15").Value = "  -1.30%  "
$ws.Range("D16").Value = "3.930.28"
$ws.Range("E16").Value = "  +0.66%  "
$ws.Range("E17").Value = "  -0.39%  "
$ws.Range("D18").Value = "19.71"
$ws.Range("E18").Value = "  -1.45%  "
$ws.Range("E19").Value = "  -2.79%  "
$ws.Range("D20").Value = "68.350.70"
$ws.Range("E20").Value = "  +1.45%  "
$ws.Range("D21").Value = "435.85"
$ws.Range("E21").Value = "  +0.41%  "
$ws.Range("D22").Value = "14.72"
$ws.Range("E22").Value = "  -1.06%  "
$ws.Range("D23").Value = "3.36"
$ws.Range("E23").Value = "  +1.44%  "
$ws.Range("D24").Value = "87.91"
$ws.Range("E24").Value = "  -0.88%  "
$ws.Range("D25").Value = "11.73"
$ws.Range("E25").Value = "  +18.44%  "
$ws.Range("D26").Value = "3.57"
$ws.Range("E26").Value = "  -0.37%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").Value = "10.46"
$ws.Range("E27").Value = "  +3.87%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "38.08"
$ws.Range("E28").Value = "  +0.74%  "
$ws.Range("D29").Value = "5.82"
$ws.Range("E29").Value = "  +4.93%  "
$ws.Range("D30").Value = "708.81"
$ws.Range("E30").Value = "  -2.58%  "
$ws.Range("E31").Value = "  -1.90%  "
$ws.Range("D32").Value = "13.32"
$ws.Range("E32").Value = "  -3.55%  "
$ws.Range("D33").Value = "2.85"
$ws.Range("E33").Value = "  +2.41%  "
$ws.Range("D34").Value = "0.0₃0929"
$ws.Range("E34").Value = "  +37.94%  "
$ws.Range("D35").Value = "41.59"
$ws.Range("E35").Value = "  -5.55%  "
$ws.Range("D36").Value = "59.31"
$ws.Range("E36").Value = "  +1.98%  "
$ws.Range("E37").Value = "  +4.49%  "
$ws.Range("E38").Value = "  -5.97%  "
$ws.Range("D39").Value = "0.998"
$ws.Range("E39").Value = "  -0.14%  "
$ws.Range("E40").Value = "  -1.83%  "
$ws.Range("D41").Value = "3.09"
$ws.Range("E41").Value = "  +11.21%  "
$ws.Range("D42").Value = "2.76"
$ws.Range("E42").Value = "  +8.44%  "
$ws.Range("D43").Value = "3.01"
$ws.Range("E43").Value = "  +3.32%  "
$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").Value = "0.141"
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").Value = "0.339"
$ws.Range("E45").Value = "  -2.20%  "
$ws.Range("E46").Value = "  -0.09%  "
$ws.Range("D47").Value = "3.42"
$ws.Range("E47").Value = "  -0.72%  "
$ws.Range("E48").Value = "  -0.15%  "
$ws.Range("D49").Value = "145.78"
$ws.Range("E49").Value = "  +1.16%  "
$ws.Range("D50").Value = "3.13"
$ws.Range("E50").Value = "  -4.53%  "
$ws.Range("D51").Value = "2.83"
$ws.Range("E51").Value = "  -2.08%  "
